# "fixed issue of yasumi day" -- refresh the 15 vocabulary-card slides
# (word / reading / meaning / page-range) with the new goi-current.pptx content.
#
# Each slide has 4 shapes: Text 0 = word, Text 1 = reading (its TextRange has
# 2 leading blank paragraphs before the actual run, so the reading is updated
# in place via Characters(Start,Length) on that 3rd paragraph rather than by
# reassigning the whole TextRange, which would collapse the blank paragraphs
# and drop the run formatting), Text 2 = meaning, Text 3 = page range.
$p = $ppt.ActivePresentation

# Slide 1: 反論 -> タンパク質
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "タンパク質"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "タンパクしつ"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "protein..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 2: 結論 -> 適応
$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "適応"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "てきおう"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "adaptation, accommodation, conformity..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 3: 賛否 -> 犯罪
$s = $p.Slides.Item(3)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "犯罪"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "はんざい"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "crime, offence, offense..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 4: 利点 -> 本質
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "本質"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "ほんしつ"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "essence, true nature, substance, reality..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 5: サポート -> まあ
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "まあ"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "just (e.g. `"`"just wait here`"`"), come now, now, now | tolerably, passably, moderately, reasonably, fairly, rather, somewha..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 6: 言い換える -> 形式
$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "形式"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "けいしき"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "form (as opposed to substance) | format, form, style, manner | formality, form | mode, form | form (bilinear, quadratic, ..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 7: レベル -> テーマ
$s = $p.Slides.Item(7)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "テーマ"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "theme, topic, subject matter, motif, project, slogan..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 8: 強調 -> 平成
$s = $p.Slides.Item(8)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "平成"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "へいせい"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "Heisei era (1989.1.8-2019.4.30)..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 9: 部分 -> 年度
$s = $p.Slides.Item(9)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "年度"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "ねんど"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "fiscal year (usu. April 1 to March 31 in Japan), financial year | academic year, school year | product year..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 10: 大 -> 形成
$s = $p.Slides.Item(10)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "形成"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "けいせい"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "formation, molding, making up, taking form, giving form to | repair (e.g. plastic surgery), replacement, -plasty..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 11: 付く -> 役に立つ
$s = $p.Slides.Item(11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "役に立つ"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "やくにたつ"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "to be helpful, to be useful..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 12: 土地 -> 確か
$s = $p.Slides.Item(12)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "確か"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "たしか"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "sure, certain, positive, definite | reliable, trustworthy, safe, sound, firm, accurate, correct, exact | If I'm not mista..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 13: 農作物 -> 学力
$s = $p.Slides.Item(13)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "学力"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "がくりょく"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "scholarly ability, scholarship, knowledge, literary ability..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 14: 矢張り -> 面
$s = $p.Slides.Item(14)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "面"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "めん"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "face | mask, face guard | (in kendo) striking the head | surface (esp. a geometrical surface) | page | aspect, facet, sid..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"

# Slide 15: 頼る -> 性
$s = $p.Slides.Item(15)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "性"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$para = $tr.Paragraphs(3)
$tr.Characters($para.Start, $para.Length).Text = "せい"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "nature (of a person) | sex, gender | sex (i.e. sexual attraction, activity, etc.) | gender | -ty, -ity, -ness, -cy..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "481-495"
